$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AllEntries")

$ws.Rows.Item(17).Insert()

Write-Host ("Row16 A: " + $ws.Range("A16").Text)
Write-Host ("Row17 A: " + $ws.Range("A17").Text)
Write-Host ("Row18 A: " + $ws.Range("A18").Text)
Write-Host ("Dimension: " + $ws.UsedRange.Address())
